# Updates the crypto price/volume table (Sheet1) to reflect the latest
# GitHub Actions scrape. For D-column price cells whose new value would
# otherwise be auto-parsed by Excel as a number (losing formatting such
# as trailing zeros or triple-dot separators), the cell's NumberFormat
# is first forced to Text ("@") so the exact string is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.938.61'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.811.30'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.24'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4985'
$ws.Range('E7').Value = '  -2.71%  '
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09756'
$ws.Range('E9').Value = '  +24.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.098'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.85'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.415'
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('B13').Value = 'BinanceUSD'
$ws.Range('C13').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.001'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.42'
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('D15').Value = '1.814.08'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.272'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001135'
$ws.Range('E17').Value = '  +5.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.17'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.16'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.906'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '27.993.10'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.08'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.245'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.74'
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('D27').Value = '2.021.30'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.51'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.379'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.36'
$ws.Range('E30').Value = '  +2.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1063'
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.548'
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.611'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06738'
$ws.Range('E35').Value = '  -4.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02324'
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.858'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2135'
$ws.Range('E38').Value = '  +0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.922'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.24'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6155'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.161'
$ws.Range('E43').Value = '  +0.82%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5878'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.288'
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.691'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.49'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.934'
$ws.Range('E49').Value = '  +1.88%  '
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06761'
$ws.Range('E51').Value = '  -1.26%  '
